{"js": "const texts = [\n  \"The four-stroke cycle\\u2014intake, compression, power, and exhaust\\u2014is the basis of most internal combustion engines. Air and fuel enter the cylinder, compressed by the piston, ignited for power, and expelled as exhaust. This process repeats rapidly to power vehicles.\",\n  \"Petrol engines use spark ignition, while diesel engines rely on compression ignition. Petrol engines are smoother and lighter, whereas diesel engines deliver torque and fuel economy, especially in commercial vehicles. Both have distinct advantages and drawbacks.\",\n  \"Forced induction technologies like turbocharging and supercharging increase engine efficiency and performance by compressing air into the cylinders. These advancements allow smaller engines to deliver higher power.\",\n  \"Electric motors differ fundamentally. They deliver instant torque, operate quietly, and rely on batteries. Hybrids combine both systems, allowing regenerative braking and improved efficiency.\"\n];\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// The document currently has 1 title paragraph followed by 5 \"chapter\"\n// paragraphs (each containing multiple runs joined with manual line\n// breaks). Replace the text of the first four chapter paragraphs with\n// the four new paragraph texts (clear() + insertText(...,\"Start\")\n// rewrites the paragraph contents as a single clean run, collapsing\n// the old multi-run / <w:br/> content).\nfor (let i = 0; i < 4; i++) {\n  const p = paras.items[i + 1];\n  p.clear();\n  p.insertText(texts[i], Word.InsertLocation.start);\n}\n\n// The fifth chapter paragraph becomes the start of the repeating block:\n// replace its text with the first paragraph text.\nlet lastPara = paras.items[5];\nlastPara.clear();\nlastPara.insertText(texts[0], Word.InsertLocation.start);\n\n// Now append the remaining 3 paragraphs of this second repeated group,\n// plus 5 more full groups of 4, for a total of 7 repeats of the\n// 4-paragraph block after the title.\nfor (let group = 0; group < 6; group++) {\n  const start = group === 0 ? 1 : 0;\n  for (let i = start; i < 4; i++) {\n    lastPara = lastPara.insertParagraph(texts[i], Word.InsertLocation.after);\n  }\n}\n\nawait context.sync();\n", "ps1": "$texts = @(\n  \"The four-stroke cycle\u2014intake, compression, power, and exhaust\u2014is the basis of most internal combustion engines. Air and fuel enter the cylinder, compressed by the piston, ignited for power, and expelled as exhaust. This process repeats rapidly to power vehicles.\",\n  \"Petrol engines use spark ignition, while diesel engines rely on compression ignition. Petrol engines are smoother and lighter, whereas diesel engines deliver torque and fuel economy, especially in commercial vehicles. Both have distinct advantages and drawbacks.\",\n  \"Forced induction technologies like turbocharging and supercharging increase engine efficiency and performance by compressing air into the cylinders. These advancements allow smaller engines to deliver higher power.\",\n  \"Electric motors differ fundamentally. They deliver instant torque, operate quietly, and rely on batteries. Hybrids combine both systems, allowing regenerative braking and improved efficiency.\"\n)\n\n$d = $word.ActiveDocument\n\n# Paragraphs 2-5 hold \"Chapter 1\"..\"Chapter 4\" (paragraph 1 is the\n# document title). Replace their text in place with the four new\n# paragraph texts; trimming the trailing paragraph mark from the range\n# before assignment collapses the old multi-run / <w:br/> content down\n# to a single run.\nfor ($i = 0; $i -lt 4; $i++) {\n    $p = $d.Paragraphs.Item($i + 2)\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.Text = $texts[$i]\n}\n\n# Paragraph 6 holds \"Chapter 5\"; it becomes the start of the repeating\n# 4-paragraph block, so replace its text with the first paragraph text.\n$p6 = $d.Paragraphs.Item(6)\n$r6 = $p6.Range\n$r6.End = $r6.End - 1\n$r6.Text = $texts[0]\n\n# Append the remaining 3 paragraphs of this second group, plus 5 more\n# full groups of 4, for a total of 7 repeats of the 4-paragraph block\n# after the title.\nfor ($group = 0; $group -lt 6; $group++) {\n    if ($group -eq 0) { $start = 1 } else { $start = 0 }\n    for ($i = $start; $i -lt 4; $i++) {\n        $endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)\n        $endRange.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n        $newPara.Range.Text = $texts[$i]\n    }\n}\n"}
